$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2154471544715447
$ws.Range("C2").Value = 0.5365853658536586
$ws.Range("J2").Value = 0.01626016260162602
$ws.Range("P2").Value = 0.1829268292682927
$ws.Range("S2").Value = 0.04878048780487805
$ws.Range("J3").Value = 0.05147058823529412
$ws.Range("P3").Value = 0.6985294117647058
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.05555555555555555
$ws.Range("P4").Value = 0.6111111111111112
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.0576923076923077
$ws.Range("D6").Value = 0.01923076923076923
$ws.Range("F6").Value = 0.08653846153846154
$ws.Range("J6").Value = 0.1442307692307692
$ws.Range("O6").Value = 0.03365384615384615
$ws.Range("Q6").Value = 0.1634615384615385
$ws.Range("R6").Value = 0.08653846153846154
$ws.Range("S6").Value = 0.4086538461538461
$ws.Range("B7").Value = 0.1118881118881119
$ws.Range("D7").Value = 0.02097902097902098
$ws.Range("F7").Value = 0.04195804195804196
$ws.Range("J7").Value = 0.09090909090909091
$ws.Range("O7").Value = 0.01398601398601399
$ws.Range("Q7").Value = 0.2167832167832168
$ws.Range("R7").Value = 0.0979020979020979
$ws.Range("S7").Value = 0.4055944055944056
$ws.Range("B8").Value = 0.09389671361502347
$ws.Range("D8").Value = 0.02347417840375587
$ws.Range("F8").Value = 0.07746478873239436
$ws.Range("J8").Value = 0.107981220657277
$ws.Range("O8").Value = 0.009389671361502348
$ws.Range("Q8").Value = 0.1854460093896714
$ws.Range("R8").Value = 0.07981220657276995
$ws.Range("S8").Value = 0.4225352112676056
$ws.Range("B9").Value = 0.084070796460177
$ws.Range("D9").Value = 0.01769911504424779
$ws.Range("F9").Value = 0.084070796460177
$ws.Range("J9").Value = 0.1106194690265487
$ws.Range("O9").Value = 0.02654867256637168
$ws.Range("Q9").Value = 0.1991150442477876
$ws.Range("R9").Value = 0.1150442477876106
$ws.Range("S9").Value = 0.3628318584070797
$ws.Range("B10").Value = 0.09193408499566348
$ws.Range("D10").Value = 0.01387684301821336
$ws.Range("F10").Value = 0.07892454466608846
$ws.Range("J10").Value = 0.1153512575888985
$ws.Range("O10").Value = 0.01474414570685169
$ws.Range("Q10").Value = 0.2298352124891587
$ws.Range("R10").Value = 0.09106678230702515
$ws.Range("S10").Value = 0.3642671292281006
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.125
$ws.Range("K11").Value = 0.2208333333333333
$ws.Range("L11").Value = 0.5083333333333333
$ws.Range("S11").Value = 0.0125
$ws.Range("G12").Value = 0.782258064516129
$ws.Range("J12").Value = 0.1370967741935484
$ws.Range("K12").Value = 0.01612903225806452
$ws.Range("L12").Value = 0.02419354838709677
$ws.Range("S12").Value = 0.04032258064516129
$ws.Range("G13").Value = 0.6333333333333333
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.03333333333333333
$ws.Range("F15").Value = 0.015
$ws.Range("H15").Value = 0.155
$ws.Range("I15").Value = 0.06
$ws.Range("J15").Value = 0.38
$ws.Range("K15").Value = 0.08500000000000001
$ws.Range("M15").Value = 0.005
$ws.Range("O15").Value = 0.05
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.0130718954248366
$ws.Range("H16").Value = 0.2156862745098039
$ws.Range("I16").Value = 0.0915032679738562
$ws.Range("J16").Value = 0.3986928104575164
$ws.Range("K16").Value = 0.08496732026143791
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("O16").Value = 0.05228758169934641
$ws.Range("S16").Value = 0.1241830065359477
$ws.Range("F17").Value = 0.006651884700665188
$ws.Range("H17").Value = 0.1951219512195122
$ws.Range("I17").Value = 0.0975609756097561
$ws.Range("J17").Value = 0.4279379157427938
$ws.Range("K17").Value = 0.09312638580931264
$ws.Range("M17").Value = 0.01108647450110865
$ws.Range("N17").Value = 0.002217294900221729
$ws.Range("O17").Value = 0.04878048780487805
$ws.Range("S17").Value = 0.1175166297117517
$ws.Range("F18").Value = 0.01015228426395939
$ws.Range("H18").Value = 0.1573604060913706
$ws.Range("I18").Value = 0.1472081218274112
$ws.Range("J18").Value = 0.4365482233502538
$ws.Range("K18").Value = 0.08121827411167512
$ws.Range("M18").Value = 0.01015228426395939
$ws.Range("N18").Value = 0.005076142131979695
$ws.Range("O18").Value = 0.06598984771573604
$ws.Range("S18").Value = 0.08629441624365482
$ws.Range("F19").Value = 0.009532062391681109
$ws.Range("H19").Value = 0.2123050259965338
$ws.Range("I19").Value = 0.1117850953206239
$ws.Range("J19").Value = 0.3760831889081456
$ws.Range("K19").Value = 0.08232235701906412
$ws.Range("M19").Value = 0.01646447140381282
$ws.Range("N19").Value = 0.0008665511265164644
$ws.Range("O19").Value = 0.07279029462738301
$ws.Range("S19").Value = 0.1178509532062392
